$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 418 (shifts existing rows 418-446 down to 419-447)
$ws.Rows.Item(418).Insert()

# Populate the newly inserted row 418 with the new weekly price record
$ws.Cells.Item(418, 1).Value = 5
$ws.Cells.Item(418, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(418, 3).Value = "Maule"
$ws.Cells.Item(418, 4).Value = 45166
$ws.Cells.Item(418, 5).Value = 7
$ws.Cells.Item(418, 6).Value = "Fruta"
$ws.Cells.Item(418, 7).Value = 100108
$ws.Cells.Item(418, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(418, 9).Value = 100108005
$ws.Cells.Item(418, 10).Value = "Piña"
$ws.Cells.Item(418, 11).Value = "Caramelo"
$ws.Cells.Item(418, 12).Value = "Segunda"
$ws.Cells.Item(418, 13).Value = 200
$ws.Cells.Item(418, 14).Value = 21000
$ws.Cells.Item(418, 15).Value = 21000
$ws.Cells.Item(418, 16).Value = 21000
$ws.Cells.Item(418, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(418, 18).Value = "Ecuador"
$ws.Cells.Item(418, 19).Value = 1500
$ws.Cells.Item(418, 20).Value = 14
